$wb = $excel.ActiveWorkbook

# "Air" sheet: cell B2 ("TO") value changes from "LAX" to "lax"
$ws = $wb.Worksheets.Item("Air")
$ws.Range("B2").Value = "lax"
